# Update the CORE_holdings workbook:
#  - refresh the "as of" date in the confidential disclaimer banner
#  - refresh the Weight / Percent Change figures for the model holdings table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no password needed to lift it from automation);
# unprotect so the cell writes below are not rejected, then restore
# protection afterwards so the sheet stays locked for end users.
$ws.Unprotect()

# --- Disclaimer banner: bump the "Model holdings provided as of" date ---
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."

# --- Row 2 (IVE / iShares S&P 500 Value ETF) ---
$ws.Range("D2").Value = 0.4927764636038424
$ws.Range("E2").Value = 0.007514174465468981

# --- Row 3 (IVW / iShares S&P Mid-Cap 400 Growth ETF) ---
$ws.Range("D3").Value = 0.2501170515853222
$ws.Range("E3").Value = -0.002587694077055946

# --- Row 4 (IJK) ---
$ws.Range("D4").Value = 0.09817873853583005
$ws.Range("E4").Value = 0.001840942562591819

# --- Row 5 (IJJ) ---
$ws.Range("D5").Value = 0.1020663047503053
$ws.Range("E5").Value = 0.01049376932446355

# --- Row 6 (IJS) ---
$ws.Range("D6").Value = 0.02938912843236561
$ws.Range("E6").Value = 0.01727840687231552

# --- Row 7 (IJT) ---
$ws.Range("D7").Value = 0.02747231309233458
$ws.Range("E7").Value = 0.004671465768111283

# --- Row 8 (Total) ---
$ws.Range("E8").Value = 0.004943516872259357

$ws.Protect()
